$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8298677802085876
$ws.Range("B1").Value = 3.012964725494385
$ws.Range("C1").Value = 3.076758623123169
$ws.Range("D1").Value = 2.606669902801514
$ws.Range("E1").Value = 2.20512843132019
